$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "code" and fill in the new "organisation" header
$ws.Columns.Item(1).Insert()
$ws.Range("A1").Value = "organisation"
$ws.Range("A1").Font.Bold = $true

# Give the new column a custom width (closest attainable value to 12.1640625)
$ws.Columns.Item(1).ColumnWidth = 11.33

# The active selection moved to C5 after the column insert
[void]$ws.Range("C5").Select()
